$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B9").Value = 4400309.12
$ws.Range("C9").Value = 707071.24
$ws.Range("D9").Value = 5107380.36
$ws.Range("E9").Value = 13.84410774528647
$ws.Range("F9").Value = 86.15589225471354
$ws.Range("G9").Value = -31.66509027373956
$ws.Range("H9").Value = -20.53649957438185
$ws.Range("I9").Value = 44196
$ws.Range("J9").Value = 1914
$ws.Range("K9").Value = 46110
$ws.Range("L9").Value = 32071
$ws.Range("M9").Value = 159.2522952199807
$ws.Range("N9").Value = 8.724785554337755
